$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "PABLO MIGUEL SAN MARTIN"
$ws.Range("C3").Value = "PABLO MIGUEL SAN MARTIN"

$ws.Range("C3").Select()
